$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the "updates automatically" date placeholders (datetimeFigureOut
#    fields) on the Slide Master, Notes Master and Handout Master from
#    11/12/19 -> 8/10/20, as PowerPoint does whenever the deck is re-saved.
# ---------------------------------------------------------------------------

# Slide Master
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "8/10/20"
    }
}

# Notes Master
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "8/10/20"
    }
}

# Handout Master
$handoutMaster = $p.HandoutMaster
for ($i = 1; $i -le $handoutMaster.Shapes.Count; $i++) {
    $shp = $handoutMaster.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "8/10/20"
    }
}

# ---------------------------------------------------------------------------
# 2) Remove the leftover "Some slides from Valentina Korzhova" credit textbox
#    (shape id 20, named "Subtitle 2") from the Title Slide custom layout.
# ---------------------------------------------------------------------------

$layout1 = $master.CustomLayouts.Item(1)
for ($i = $layout1.Shapes.Count; $i -ge 1; $i--) {
    $shp = $layout1.Shapes.Item($i)
    if ($shp.Id -eq 20 -and $shp.Name -eq "Subtitle 2") {
        $shp.Delete()
    }
}
